$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 1.68
$ws.Cells.Item(2, 14).Value = 4.1
$ws.Cells.Item(2, 28).Value = 8.6
$ws.Cells.Item(3, 6).Value = 1.09
$ws.Cells.Item(3, 7).Value = 40
$ws.Cells.Item(3, 8).Value = 1.13
$ws.Cells.Item(3, 9).Value = 1.18
$ws.Cells.Item(3, 10).Value = 6.6
$ws.Cells.Item(3, 11).Value = 11
$ws.Cells.Item(3, 12).Value = 1.22
$ws.Cells.Item(3, 14).Value = 4.1
$ws.Cells.Item(3, 15).Value = 1.18
$ws.Cells.Item(3, 16).Value = 2.4
$ws.Cells.Item(3, 17).Value = 1.53
$ws.Cells.Item(3, 18).Value = 1.56
$ws.Cells.Item(3, 19).Value = 2.18
$ws.Cells.Item(3, 20).Value = 2.72
$ws.Cells.Item(3, 21).Value = 1.43
$ws.Cells.Item(3, 22).Value = 5.7
$ws.Cells.Item(3, 23).Value = 1.02
$ws.Cells.Item(3, 24).Value = 34
$ws.Cells.Item(3, 25).Value = 10.5
$ws.Cells.Item(3, 26).Value = 8.199999999999999
$ws.Cells.Item(3, 27).Value = 8.800000000000001
$ws.Cells.Item(3, 28).Value = 95
$ws.Cells.Item(3, 29).Value = 27
$ws.Cells.Item(3, 30).Value = 17.5
$ws.Cells.Item(3, 31).Value = 21
$ws.Cells.Item(3, 33).Value = 990
$ws.Cells.Item(3, 34).Value = 85
$ws.Cells.Item(3, 35).Value = 90
$ws.Cells.Item(3, 41).Value = 4.2
$ws.Cells.Item(4, 6).Value = 2.48
$ws.Cells.Item(4, 7).Value = 3.25
$ws.Cells.Item(4, 8).Value = 2.72
$ws.Cells.Item(4, 9).Value = 3.6
$ws.Cells.Item(4, 10).Value = 2.8
$ws.Cells.Item(4, 11).Value = 4.2
$ws.Cells.Item(4, 12).Value = 1.39
$ws.Cells.Item(4, 13).Value = 1.07
$ws.Cells.Item(4, 14).Value = 2.78
$ws.Cells.Item(4, 15).Value = 1.37
$ws.Cells.Item(4, 16).Value = 1.69
$ws.Cells.Item(4, 17).Value = 1.96
$ws.Cells.Item(4, 18).Value = 1.26
$ws.Cells.Item(4, 19).Value = 3.45
$ws.Cells.Item(4, 20).Value = 1.8
$ws.Cells.Item(4, 21).Value = 1.95
$ws.Cells.Item(4, 22).Value = 1.39
$ws.Cells.Item(4, 23).Value = 1.45
$ws.Cells.Item(5, 6).Value = 5.5
$ws.Cells.Item(5, 7).Value = 6.6
$ws.Cells.Item(5, 8).Value = 1.6
$ws.Cells.Item(5, 9).Value = 1.67
$ws.Cells.Item(5, 12).Value = 1.28
$ws.Cells.Item(5, 16).Value = 2.12
$ws.Cells.Item(5, 18).Value = 1.44
$ws.Cells.Item(5, 20).Value = 1.79
$ws.Cells.Item(5, 22).Value = 2.46
$ws.Cells.Item(5, 23).Value = 1.18
$ws.Cells.Item(5, 25).Value = 10
$ws.Cells.Item(5, 27).Value = 17
$ws.Cells.Item(5, 30).Value = 10.5
$ws.Cells.Item(5, 32).Value = 48
$ws.Cells.Item(5, 34).Value = 20
$ws.Cells.Item(5, 35).Value = 32
$ws.Cells.Item(5, 36).Value = 160
$ws.Cells.Item(5, 37).Value = 90
$ws.Cells.Item(5, 39).Value = 100
$ws.Cells.Item(5, 40).Value = 85
$ws.Cells.Item(5, 41).Value = 8.800000000000001
$ws.Cells.Item(6, 6).Value = 1.41
$ws.Cells.Item(6, 7).Value = 1.47
$ws.Cells.Item(6, 8).Value = 8.4
$ws.Cells.Item(6, 10).Value = 4.7
$ws.Cells.Item(6, 12).Value = 1.28
$ws.Cells.Item(6, 15).Value = 1.24
$ws.Cells.Item(6, 16).Value = 2.16
$ws.Cells.Item(6, 18).Value = 1.45
$ws.Cells.Item(6, 19).Value = 2.8
$ws.Cells.Item(6, 20).Value = 1.97
$ws.Cells.Item(6, 21).Value = 1.87
$ws.Cells.Item(6, 22).Value = 1.1
$ws.Cells.Item(6, 23).Value = 3.05
$ws.Cells.Item(6, 37).Value = 18.5
$ws.Cells.Item(6, 40).Value = 7.8
$ws.Cells.Item(7, 6).Value = 1.83
$ws.Cells.Item(7, 7).Value = 1.86
$ws.Cells.Item(7, 12).Value = 1.34
$ws.Cells.Item(7, 14).Value = 3.6
$ws.Cells.Item(7, 16).Value = 1.9
$ws.Cells.Item(7, 17).Value = 1.94
$ws.Cells.Item(7, 19).Value = 3.45
$ws.Cells.Item(7, 23).Value = 2.16
$ws.Cells.Item(7, 24).Value = 18
$ws.Cells.Item(8, 7).Value = 1.43
$ws.Cells.Item(8, 14).Value = 4
$ws.Cells.Item(8, 17).Value = 1.93
$ws.Cells.Item(8, 18).Value = 1.4
$ws.Cells.Item(8, 21).Value = 1.79
$ws.Cells.Item(8, 23).Value = 3.3
$ws.Cells.Item(8, 26).Value = 90
$ws.Cells.Item(8, 28).Value = 7.6
$ws.Cells.Item(8, 38).Value = 42
$ws.Cells.Item(8, 39).Value = 200
$ws.Cells.Item(10, 6).Value = 3.85
$ws.Cells.Item(10, 7).Value = 4.3
$ws.Cells.Item(10, 9).Value = 2.2
$ws.Cells.Item(10, 10).Value = 3.45
$ws.Cells.Item(10, 17).Value = 1.96
$ws.Cells.Item(11, 8).Value = 1.62
$ws.Cells.Item(11, 14).Value = 3.4
$ws.Cells.Item(11, 21).Value = 1.78
$ws.Cells.Item(11, 35).Value = 44
$ws.Cells.Item(12, 8).Value = 2.9
$ws.Cells.Item(12, 21).Value = 1.98
$ws.Cells.Item(14, 6).Value = 3.5
$ws.Cells.Item(14, 7).Value = 3.75
$ws.Cells.Item(14, 8).Value = 2.42
$ws.Cells.Item(14, 11).Value = 3.15
$ws.Cells.Item(14, 16).Value = 1.58
$ws.Cells.Item(14, 22).Value = 1.65
$ws.Cells.Item(14, 23).Value = 1.37
$ws.Cells.Item(15, 7).Value = 2.16
$ws.Cells.Item(15, 19).Value = 2.94
$ws.Cells.Item(15, 20).Value = 1.67
$ws.Cells.Item(15, 22).Value = 1.31
$ws.Cells.Item(15, 23).Value = 1.86
$ws.Cells.Item(16, 6).Value = 2.28
$ws.Cells.Item(16, 7).Value = 2.32
$ws.Cells.Item(16, 21).Value = 2.16
$ws.Cells.Item(16, 23).Value = 1.76
$ws.Cells.Item(16, 26).Value = 24
$ws.Cells.Item(16, 32).Value = 13.5
$ws.Cells.Item(16, 40).Value = 18.5
$ws.Cells.Item(17, 6).Value = 1.8
$ws.Cells.Item(17, 7).Value = 1.82
$ws.Cells.Item(17, 8).Value = 5.5
$ws.Cells.Item(17, 10).Value = 3.75
$ws.Cells.Item(17, 11).Value = 3.85
$ws.Cells.Item(17, 18).Value = 1.34
$ws.Cells.Item(17, 21).Value = 1.99
$ws.Cells.Item(17, 23).Value = 2.22
$ws.Cells.Item(18, 9).Value = 2.22
$ws.Cells.Item(18, 15).Value = 1.38
$ws.Cells.Item(19, 7).Value = 4
$ws.Cells.Item(19, 8).Value = 2.1
$ws.Cells.Item(19, 11).Value = 3.7
$ws.Cells.Item(19, 14).Value = 4.5
$ws.Cells.Item(19, 24).Value = 16.5
$ws.Cells.Item(20, 16).Value = 3.1
$ws.Cells.Item(20, 17).Value = 1.42
$ws.Cells.Item(20, 18).Value = 1.86
$ws.Cells.Item(20, 19).Value = 2.04
$ws.Cells.Item(20, 20).Value = 2.32
$ws.Cells.Item(20, 21).Value = 1.69
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(8, 11).Value = 5.1

Write-Output "Updated 158 cells"
